$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header figures: total mora value and worker count went up (new worker
#    added to the account statement).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 235251
$ws.Range("C13").Value = 3

# ---------------------------------------------------------------------------
# 2. Make room for 3 extra data rows (table grows from 6 to 9 rows) by
#    inserting before the old row 22 (i.e. right after the last existing
#    data row, 21). This pushes the trailing signature block (old rows
#    26/27) down to rows 29/30 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("22:24").Insert()

# Copy the "interior" row formatting down into the 2 new interior rows.
$ws.Range("B20:J20").Copy()
$ws.Range("B22:J23").PasteSpecial(-4122)

# Copy the "last row" (bottom border) formatting into the new last row.
$ws.Range("B21:J21").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

# Old row 21 is no longer the last table row, so it needs to switch from the
# bottom-border style back to the regular interior style.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Rewrite the worker detail table (rows 16-24): a new worker (JOSE MANUEL
#    RODRIGUEZ SERRANO) is added, and the period ordering per worker is
#    normalized to 2004, 2003, 2002 for each of the 3 workers in turn.
# ---------------------------------------------------------------------------
$docType = "CC"
$period2004 = "2004"
$period2003 = "2003"
$period2002 = "2002"
$valor2004 = 8193
$valor2003 = 35112
$valor2002 = 35112
$salario = 877803

$workers = @(
    @{ Doc = "1047372540"; Name = "JOSE MANUEL RODRIGUEZ SERRANO" },
    @{ Doc = "37285383";   Name = "MYRIAM SALCEDO SALCEDO" },
    @{ Doc = "1004823600"; Name = "JOSE IVAN PEREZ SALCEDO" }
)

$row = 16
foreach ($worker in $workers) {
    $ws.Cells.Item($row, 2).Value = $docType
    $ws.Cells.Item($row, 3).Value = $worker.Doc
    $ws.Cells.Item($row, 4).Value = $worker.Name
    $ws.Cells.Item($row, 5).Value = $period2004
    $ws.Cells.Item($row, 6).Value = $valor2004
    $ws.Cells.Item($row, 7).Value = $salario
    $row = $row + 1

    $ws.Cells.Item($row, 2).Value = $docType
    $ws.Cells.Item($row, 3).Value = $worker.Doc
    $ws.Cells.Item($row, 4).Value = $worker.Name
    $ws.Cells.Item($row, 5).Value = $period2003
    $ws.Cells.Item($row, 6).Value = $valor2003
    $ws.Cells.Item($row, 7).Value = $salario
    $row = $row + 1

    $ws.Cells.Item($row, 2).Value = $docType
    $ws.Cells.Item($row, 3).Value = $worker.Doc
    $ws.Cells.Item($row, 4).Value = $worker.Name
    $ws.Cells.Item($row, 5).Value = $period2002
    $ws.Cells.Item($row, 6).Value = $valor2002
    $ws.Cells.Item($row, 7).Value = $salario
    $row = $row + 1
}
